$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Rv1165", 3, "typA bipA Rv1165", "FUNCTION: A 50S ribosomal subunit assembly protein with GTPase activity, required for 50S subunit assembly at low temperatures, may also play a role in translation. Binds GTP and analogs. Binds the 70S ribosome between the 30S and 50S subunits, in a similar position as ribosome-bound EF-G; it contacts a number of ribosomal proteins, both rRNAs and the A-site tRNA. {ECO:0000256|HAMAP-Rule:MF_00849}.", 27),
    @("Rv0319", 2, "pcp Rv0319 MTCY63.24", "FUNCTION: Removes 5-oxoproline from various penultimate amino acid residues except L-proline. {ECO:0000250}.", 27),
    @("Rv0795", 1, "Rv0795 MTV042.05; Rv1370c MTCY02B12.04c; Rv1757c MTCY28.23c; Rv1763 MTCY28.29; Rv2105 MTCY261.01; Rv2168c MTV021.01c; Rv2278 MTCY339.32c; Rv2354 MTCY98.23; Rv2480c MTV008.36c; Rv2648 MTCY441.17A; Rv2815c MTCY16B7.28; Rv3184 MYV014.28; Rv3186 MTV014.30; Rv3325 MTV016.25; Rv3381c MTV004.39c; Rv3474 MTCY13E12.27", "", 27),
    @("Rv0218", 1, "Rv0218", "", 27),
    @("Rv0336", 1, "Rv0336", "", 27),
    @("Rv0325", 1, "Rv0325", "", 27),
    @("Rv0399c", 1, "lpqK Rv0399c", "", 27),
    @("Rv0585c", 1, "Rv0585c", "", 27),
    @("Rv1768", 1, "PE_PGRS31 Rv1768", "", 27),
    @("Rv0976c", 1, "Rv0976c", "", 27),
    @("Rv3629c", 1, "Rv3629c", "", 27)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
